# Apply the "corrected TwelveBarBreakout" results to the S32 workbook:
#   - Signal Instances: append 12 new twelve_bar_breakout signal rows (53-64)
#   - Summary: insert a twelve_bar_breakout aggregate row ahead of random_baseline

$wb = $excel.ActiveWorkbook
$wsSignals = $wb.Worksheets.Item("Signal Instances")
$wsSummary = $wb.Worksheets.Item("Summary")

# --- Signal Instances: 12 new twelve_bar_breakout rows (columns A:P), one per signal_bar ---
# Columns: signal_type, signal_bar, signal_date, entry_price, return_12m, profitable_12m,
#          mfe_12m, mfe_bar, mae_12m, mae_bar, exit_signal_fired, exit_signal_bar,
#          return_at_exit, left_on_table, exit_vs_hold, exit_vs_mfe
$sheet1Data = New-Object 'object[,]' 12,16
# Row 53: signal_bar=48
$sheet1Data[0,0] = "twelve_bar_breakout"
$sheet1Data[0,1] = 48
$sheet1Data[0,2] = ""
$sheet1Data[0,3] = 1.172122478485107
$sheet1Data[0,4] = 70.38532308647385
$sheet1Data[0,5] = $true
$sheet1Data[0,6] = 84.65614545788821
$sheet1Data[0,7] = 33
$sheet1Data[0,8] = -13.65268929941376
$sheet1Data[0,9] = 6
$sheet1Data[0,10] = $true
$sheet1Data[0,11] = 40
$sheet1Data[0,12] = 71.70309899561413
$sheet1Data[0,13] = 14.27082237141437
$sheet1Data[0,14] = 1.317775909140281
$sheet1Data[0,15] = -12.95304646227409
# Row 54: signal_bar=116
$sheet1Data[1,0] = "twelve_bar_breakout"
$sheet1Data[1,1] = 116
$sheet1Data[1,2] = ""
$sheet1Data[1,3] = 2.14080286026001
$sheet1Data[1,4] = 28.49539233992243
$sheet1Data[1,5] = $true
$sheet1Data[1,6] = 44.0705789265018
$sheet1Data[1,7] = 40
$sheet1Data[1,8] = -5.369143571519596
$sheet1Data[1,9] = 3
$sheet1Data[1,10] = $true
$sheet1Data[1,11] = 18
$sheet1Data[1,12] = 10.79237988154374
$sheet1Data[1,13] = 15.57518658657936
$sheet1Data[1,14] = -17.70301245837869
$sheet1Data[1,15] = -33.27819904495805
# Row 55: signal_bar=121
$sheet1Data[2,0] = "twelve_bar_breakout"
$sheet1Data[2,1] = 121
$sheet1Data[2,2] = ""
$sheet1Data[2,3] = 2.234192848205566
$sheet1Data[2,4] = 27.1944567876414
$sheet1Data[2,5] = $true
$sheet1Data[2,6] = 38.04838185427421
$sheet1Data[2,7] = 35
$sheet1Data[2,8] = -1.607715562792441
$sheet1Data[2,9] = 0
$sheet1Data[2,10] = $true
$sheet1Data[2,11] = 13
$sheet1Data[2,12] = 6.161222356396705
$sheet1Data[2,13] = 10.85392506663282
$sheet1Data[2,14] = -21.03323443124469
$sheet1Data[2,15] = -31.88715949787751
# Row 56: signal_bar=156
$sheet1Data[3,0] = "twelve_bar_breakout"
$sheet1Data[3,1] = 156
$sheet1Data[3,2] = ""
$sheet1Data[3,3] = 3.038798809051514
$sheet1Data[3,4] = -7.808274292256352
$sheet1Data[3,5] = $false
$sheet1Data[3,6] = 9.333800452459423
$sheet1Data[3,7] = 20
$sheet1Data[3,8] = -23.36416610507271
$sheet1Data[3,9] = 30
$sheet1Data[3,10] = $true
$sheet1Data[3,11] = 5
$sheet1Data[3,12] = -9.476302731144424
$sheet1Data[3,13] = 17.14207474471577
$sheet1Data[3,14] = -1.668028438888071
$sheet1Data[3,15] = -18.81010318360385
# Row 57: signal_bar=199
$sheet1Data[4,0] = "twelve_bar_breakout"
$sheet1Data[4,1] = 199
$sheet1Data[4,2] = ""
$sheet1Data[4,3] = 3.032327890396118
$sheet1Data[4,4] = -47.68670924766708
$sheet1Data[4,5] = $false
$sheet1Data[4,6] = 2.887159828418477
$sheet1Data[4,7] = 3
$sheet1Data[4,8] = -56.20282492283049
$sheet1Data[4,9] = 52
$sheet1Data[4,10] = $true
$sheet1Data[4,11] = 12
$sheet1Data[4,12] = -14.96061840601549
$sheet1Data[4,13] = 50.57386907608556
$sheet1Data[4,14] = 32.7260908416516
$sheet1Data[4,15] = -17.84777823443396
# Row 58: signal_bar=238
$sheet1Data[5,0] = "twelve_bar_breakout"
$sheet1Data[5,1] = 238
$sheet1Data[5,2] = ""
$sheet1Data[5,3] = 2.326508045196533
$sheet1Data[5,4] = -8.859855861018589
$sheet1Data[5,5] = $false
$sheet1Data[5,6] = 2.083350369467065
$sheet1Data[5,7] = 1
$sheet1Data[5,8] = -44.14888023906098
$sheet1Data[5,9] = 15
$sheet1Data[5,10] = $true
$sheet1Data[5,11] = 7
$sheet1Data[5,12] = -8.333328209375573
$sheet1Data[5,13] = 10.94320623048565
$sheet1Data[5,14] = 0.5265276516430166
$sheet1Data[5,15] = -10.41667857884264
# Row 59: signal_bar=287
$sheet1Data[6,0] = "twelve_bar_breakout"
$sheet1Data[6,1] = 287
$sheet1Data[6,2] = ""
$sheet1Data[6,3] = 1.914120078086853
$sheet1Data[6,4] = 55.75096736255654
$sheet1Data[6,5] = $true
$sheet1Data[6,6] = 80.60011509413567
$sheet1Data[6,7] = 47
$sheet1Data[6,8] = -4.310340840949636
$sheet1Data[6,9] = 0
$sheet1Data[6,10] = $true
$sheet1Data[6,11] = 35
$sheet1Data[6,12] = 27.99399731067262
$sheet1Data[6,13] = 24.84914773157913
$sheet1Data[6,14] = -27.75697005188393
$sheet1Data[6,15] = -52.60611778346306
# Row 60: signal_bar=296
$sheet1Data[7,0] = "twelve_bar_breakout"
$sheet1Data[7,1] = 296
$sheet1Data[7,2] = ""
$sheet1Data[7,3] = 2.285393238067627
$sheet1Data[7,4] = 52.00406149452399
$sheet1Data[7,5] = $true
$sheet1Data[7,6] = 59.06539543652482
$sheet1Data[7,7] = 51
$sheet1Data[7,8] = -10.83029672136926
$sheet1Data[7,9] = 2
$sheet1Data[7,10] = $true
$sheet1Data[7,11] = 26
$sheet1Data[7,12] = 7.200754796187711
$sheet1Data[7,13] = 7.061333942000829
$sheet1Data[7,14] = -44.80330669833628
$sheet1Data[7,15] = -51.86464064033711
# Row 61: signal_bar=345
$sheet1Data[8,0] = "twelve_bar_breakout"
$sheet1Data[8,1] = 345
$sheet1Data[8,2] = ""
$sheet1Data[8,3] = 3.405941963195801
$sheet1Data[8,4] = 8.542579651782543
$sheet1Data[8,5] = $true
$sheet1Data[8,6] = 37.75259595576605
$sheet1Data[8,7] = 14
$sheet1Data[8,8] = -13.10597571572787
$sheet1Data[8,9] = 29
$sheet1Data[8,10] = $true
$sheet1Data[8,11] = 24
$sheet1Data[8,12] = 6.317394996736564
$sheet1Data[8,13] = 29.21001630398351
$sheet1Data[8,14] = -2.225184655045979
$sheet1Data[8,15] = -31.43520095902949
# Row 62: signal_bar=391
$sheet1Data[9,0] = "twelve_bar_breakout"
$sheet1Data[9,1] = 391
$sheet1Data[9,2] = ""
$sheet1Data[9,3] = 3.613717317581177
$sheet1Data[9,4] = -16.31245147053555
$sheet1Data[9,5] = $false
$sheet1Data[9,6] = 24.55242979699103
$sheet1Data[9,7] = 10
$sheet1Data[9,8] = -19.21643925923511
$sheet1Data[9,9] = 52
$sheet1Data[9,10] = $true
$sheet1Data[9,11] = 27
$sheet1Data[9,12] = 2.165747654108203
$sheet1Data[9,13] = 40.86488126752658
$sheet1Data[9,14] = 18.47819912464375
$sheet1Data[9,15] = -22.38668214288283
# Row 63: signal_bar=469
$sheet1Data[10,0] = "twelve_bar_breakout"
$sheet1Data[10,1] = 469
$sheet1Data[10,2] = ""
$sheet1Data[10,3] = 3.594549655914306
$sheet1Data[10,4] = -18.88226509975634
$sheet1Data[10,5] = $false
$sheet1Data[10,6] = 6.914891177257577
$sheet1Data[10,7] = 3
$sheet1Data[10,8] = -32.31053581011737
$sheet1Data[10,9] = 47
$sheet1Data[10,10] = $true
$sheet1Data[10,11] = 9
$sheet1Data[10,12] = -9.042548887238347
$sheet1Data[10,13] = 25.79715627701392
$sheet1Data[10,14] = 9.839716212517992
$sheet1Data[10,15] = -15.95744006449592
# Row 64: signal_bar=494
$sheet1Data[11,0] = "twelve_bar_breakout"
$sheet1Data[11,1] = 494
$sheet1Data[11,2] = ""
$sheet1Data[11,3] = 3.685718536376953
$sheet1Data[11,4] = -16.97684148753505
$sheet1Data[11,5] = $false
$sheet1Data[11,6] = 3.815789019632208
$sheet1Data[11,7] = 0
$sheet1Data[11,8] = -33.9848830529685
$sheet1Data[11,9] = 22
$sheet1Data[11,10] = $true
$sheet1Data[11,11] = 5
$sheet1Data[11,12] = -10.52631136351169
$sheet1Data[11,13] = 20.79263050716726
$sheet1Data[11,14] = 6.450530124023354
$sheet1Data[11,15] = -14.3421003831439

$wsSignals.Range("A53:P64").Value = $sheet1Data

# --- Summary: insert a row at 4 (shifts existing random_baseline row down to row 5), ---
# --- then fill it in with the twelve_bar_breakout aggregate stats (columns A:M; N:T stay blank) ---
$wsSummary.Rows.Item(4).Insert()

# Columns: signal_type, total_signals, win_rate_12m, mean_return_12m, median_return_12m,
#          std_return_12m, mean_mfe_12m, mean_mae_12m, mean_left_on_table, exit_fired_rate,
#          mean_exit_bar, exit_useful_rate, mean_exit_vs_hold, total_samples, baseline_win_rate,
#          baseline_mean_return, baseline_median_return, baseline_std_return, baseline_mean_mfe,
#          baseline_mean_mae
$sheet2Row4 = New-Object 'object[,]' 1,20
$sheet2Row4[0,0] = "twelve_bar_breakout"
$sheet2Row4[0,1] = 12
$sheet2Row4[0,2] = 50
$sheet2Row4[0,3] = 10.48719860534432
$sheet2Row4[0,4] = 0.3671526797630951
$sheet2Row4[0,5] = 36.17602116891146
$sheet2Row4[0,6] = 32.81505278077638
$sheet2Row4[0,7] = -21.50865759175481
$sheet2Row4[0,8] = 22.32785417543207
$sheet2Row4[0,9] = 100
$sheet2Row4[0,10] = 18.41666666666667
$sheet2Row4[0,11] = 50
$sheet2Row4[0,12] = -3.820908072513136
$sheet2Row4[0,13] = ""
$sheet2Row4[0,14] = ""
$sheet2Row4[0,15] = ""
$sheet2Row4[0,16] = ""
$sheet2Row4[0,17] = ""
$sheet2Row4[0,18] = ""
$sheet2Row4[0,19] = ""

$wsSummary.Range("A4:T4").Value = $sheet2Row4

Write-Output "Edit complete: Signal Instances A53:P64, Summary row 4 inserted"
